$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (F column) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 102
$wsExpo.Range("F3").Value = 2085
$wsExpo.Range("F4").Value = 861
$wsExpo.Range("F5").Value = 1244

# Sheet "全部类型" (all types) - same updates mirrored on the consolidated sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 102
$wsAll.Range("F3").Value = 2085
$wsAll.Range("F6").Value = 861
$wsAll.Range("F7").Value = 1244
